$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.641.32"
$ws.Range("E2").Value = "  +1.16%  "

$ws.Range("D3").Value = "1.867.87"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("E4").Value = "  +0.48%  "

$ws.Range("D5").Value = "'331.82"
$ws.Range("E5").Value = "  +2.94%  "

$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("D7").Value = "'0.4690"
$ws.Range("E7").Value = "  +3.86%  "

$ws.Range("D8").Value = "'0.3936"
$ws.Range("E8").Value = "  +2.13%  "

$ws.Range("D9").Value = "'47.88"
$ws.Range("E9").Value = "  -0.26%  "

$ws.Range("D10").Value = "'0.08055"
$ws.Range("E10").Value = "  +2.16%  "

$ws.Range("D11").Value = "'1.021"
$ws.Range("E11").Value = "  -0.19%  "

$ws.Range("D12").Value = "'21.76"
$ws.Range("E12").Value = "  +1.86%  "

$ws.Range("D13").Value = "1.865.47"
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").Value = "'5.933"
$ws.Range("E14").Value = "  +0.99%  "

$ws.Range("D15").Value = "'7.133"
$ws.Range("E15").Value = "  -0.52%  "

$ws.Range("E16").Value = "  +0.48%  "

$ws.Range("D17").Value = "'0.00001046"
$ws.Range("E17").Value = "  +1.52%  "

$ws.Range("D18").Value = "'86.61"
$ws.Range("E18").Value = "  +1.46%  "

$ws.Range("D19").Value = "'0.06627"
$ws.Range("E19").Value = "  +1.59%  "

$ws.Range("D20").Value = "'17.20"
$ws.Range("E20").Value = "  +0.96%  "

$ws.Range("E21").Value = "  +0.29%  "

$ws.Range("D22").Value = "27.663.05"
$ws.Range("E22").Value = "  +1.25%  "

$ws.Range("E23").Value = "  -0.42%  "

$ws.Range("D24").Value = "'10.99"
$ws.Range("E24").Value = "  +2.14%  "

$ws.Range("D25").Value = "'2.309"
$ws.Range("E25").Value = "  +1.92%  "

$ws.Range("D26").Value = "2.093.92"
$ws.Range("E26").Value = "  +0.73%  "

$ws.Range("D27").Value = "'158.69"
$ws.Range("E27").Value = "  +4.55%  "

$ws.Range("D28").Value = "'20.20"
$ws.Range("E28").Value = "  +2.45%  "

$ws.Range("D29").Value = "'2.088"
$ws.Range("E29").Value = "  +1.18%  "

$ws.Range("D30").Value = "'5.546"
$ws.Range("E30").Value = "  +1.26%  "

$ws.Range("D31").Value = "'122.21"
$ws.Range("E31").Value = "  +1.45%  "

$ws.Range("D32").Value = "'0.9643"
$ws.Range("E32").Value = "  +3.10%  "

$ws.Range("D33").Value = "'0.09484"
$ws.Range("E33").Value = "  +1.98%  "

$ws.Range("E34").Value = "  -2.74%  "

$ws.Range("E35").Value = "  -0.20%  "

$ws.Range("D36").Value = "'5.312"
$ws.Range("E36").Value = "  +0.59%  "

$ws.Range("D37").Value = "'0.02257"
$ws.Range("E37").Value = "  +1.25%  "

$ws.Range("D38").Value = "'0.06082"
$ws.Range("E38").Value = "  +1.69%  "

$ws.Range("D39").Value = "'1.228"
$ws.Range("E39").Value = "  +1.43%  "

$ws.Range("D40").Value = "'8.109"
$ws.Range("E40").Value = "  -2.07%  "

$ws.Range("E41").Value = "  +0.32%  "

$ws.Range("D42").Value = "'0.5980"
$ws.Range("E42").Value = "  +1.27%  "

$ws.Range("D43").Value = "'0.1892"
$ws.Range("E43").Value = "  +0.29%  "

$ws.Range("D44").Value = "'10.21"
$ws.Range("E44").Value = "  +0.62%  "

$ws.Range("E45").Value = "  -0.75%  "

$ws.Range("D46").Value = "'0.5701"
$ws.Range("E46").Value = "  +1.14%  "

$ws.Range("D47").Value = "'12.23"
$ws.Range("E47").Value = "  +2.80%  "

$ws.Range("E48").Value = "  +1.04%  "

$ws.Range("D49").Value = "'1.933"
$ws.Range("E49").Value = "  +0.46%  "

$ws.Range("D50").Value = "'0.06851"
$ws.Range("E50").Value = "  +0.72%  "

$ws.Range("D51").Value = "'114.39"
$ws.Range("E51").Value = "  +5.77%  "
